$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so values like "1.001"
# or "23.531.58" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Swap rows 8 and 9 (Cardano <-> OKB) ---
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "51.20"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.3600"
$ws.Range("E9").Value = "  +2.45%  "

# --- Update D/E price + volume values ---
$ws.Range("D2").Value = "23.531.58"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "1.655.43"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "302.26"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.3837"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D10").Value = "0.08205"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("D11").Value = "1.239"
$ws.Range("E11").Value = "  +3.99%  "
$ws.Range("D12").Value = "0.9991"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "22.41"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "6.485"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "7.503"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("D16").Value = "0.00001226"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "1.646.12"
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("D18").Value = "97.53"
$ws.Range("E18").Value = "  +3.69%  "
$ws.Range("D19").Value = "0.06990"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "6.798"
$ws.Range("E20").Value = "  +5.10%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D24").Value = "23.553.30"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "2.520"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "3.030"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "21.25"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").Value = "152.58"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "5.240"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "134.23"
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("D31").Value = "1.835.97"
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("D32").Value = "7.134"
$ws.Range("E32").Value = "  +11.92%  "
$ws.Range("D33").Value = "2.248"
$ws.Range("E33").Value = "  +6.16%  "
$ws.Range("D34").Value = "12.02"
$ws.Range("E34").Value = "  +5.84%  "
$ws.Range("D35").Value = "1.063"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "0.02799"
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("D37").Value = "6.109"
$ws.Range("E37").Value = "  +5.27%  "
$ws.Range("D38").Value = "0.2498"
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("D39").Value = "0.08766"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").Value = "0.07004"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").Value = "13.17"
$ws.Range("E41").Value = "  +10.68%  "
$ws.Range("D42").Value = "0.6995"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("D43").Value = "1.332"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "15.93"
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("D45").Value = "0.6538"
$ws.Range("E45").Value = "  +4.50%  "
$ws.Range("D46").Value = "0.9995"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").Value = "2.306"
$ws.Range("E47").Value = "  +3.28%  "
$ws.Range("D48").Value = "3.956"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "0.07899"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").Value = "128.30"
$ws.Range("E50").Value = "  +0.53%  "

# --- Update E-only rows (volume changed, price text unchanged) ---
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("E23").Value = "  +3.87%  "
$ws.Range("E51").Value = "  +2.57%  "

# Restore default style on column D now that text values are set,
# so no lingering explicit number-format style remains on the cells.
$ws.Range("D2:D51").Style = "Normal"
